$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("To Do")
$row = $ws.Rows.Item(139)
$row.UseStandardHeight = $true
$row.RowHeight = 30
Write-Host "done"
